# Generate Report for Handoff
# Adds a new localization-status row for the file
#   c61b92bc-90f8-4cb8-b3b3-410cc0e93bd0.md
# to every worksheet (Overview, zh-cn, de-de), mirroring the existing
# 166d75bf-414f-42d8-b951-2059586c0f4d.md row.

$wb = $excel.ActiveWorkbook

$newGuid   = "c61b92bc-90f8-4cb8-b3b3-410cc0e93bd0"
$newFile   = "$newGuid.md"
$newPath   = "e2e\$newGuid.md"
$commit    = "942927156999938cf922229da6e833ccb0bedfdc"
$ghUrl     = "https://github.com/OpenLocalizationTestOrg/oltest/blob/$commit/e2e/$newFile"
$xlfHash   = "971feed1a12a0c758c942732baf47a45590b6511"
$dateFmt   = "yyyy-mm-dd HH:mm:ss"

$hyperlinkColor = 15570276  # matches the workbook's existing HyperLink font color (FF6495ED)

function Style-DateCell($range) {
    $range.NumberFormat = $dateFmt
}

function Style-HyperlinkCell($range) {
    $range.Font.Underline = $true
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1) -> table3 "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A3").Value = $newFile
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-13 22:52:38"
Style-DateCell $wsOverview.Range("G3")

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $ghUrl, "", "", $newPath) | Out-Null
Style-HyperlinkCell $wsOverview.Range("B3")

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) -> table1 "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = "$newGuid.$xlfHash.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-13 22:52:30"
Style-DateCell $wsZhCn.Range("H3")
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
Style-DateCell $wsZhCn.Range("K3")
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("O3").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $ghUrl, "", "", $newFile) | Out-Null
Style-HyperlinkCell $wsZhCn.Range("A3")

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3) -> table2 "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = "$newGuid.$xlfHash.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-13 22:52:38"
Style-DateCell $wsDeDe.Range("H3")
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
Style-DateCell $wsDeDe.Range("K3")
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("O3").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $ghUrl, "", "", $newFile) | Out-Null
Style-HyperlinkCell $wsDeDe.Range("A3")

Write-Output "Added handoff row for $newFile to Overview, zh-cn and de-de sheets."
